$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing fatality counts for the last 16 day-columns (JA2:JP2)
$updates = @{
    261 = 259   # JA2
    262 = 264   # JB2
    263 = 266   # JC2
    264 = 268   # JD2
    265 = 272   # JE2
    266 = 272   # JF2
    267 = 277   # JG2
    268 = 279   # JH2
    269 = 280   # JI2
    270 = 284   # JJ2
    271 = 285   # JK2
    272 = 285   # JL2
    273 = 285   # JM2
    274 = 286   # JN2
    275 = 287   # JO2
    276 = 288   # JP2
}

foreach ($col in $updates.Keys) {
    $ws.Cells.Item(2, $col).Value = $updates[$col]
}

# Add two new day columns: JQ (277) "Fatalities 12-07" and JR (278) "Fatalities 12-08"
$ws.Cells.Item(1, 277).Value = "Fatalities 12-07"
$ws.Cells.Item(1, 278).Value = "Fatalities 12-08"
$ws.Cells.Item(2, 277).Value = 288
$ws.Cells.Item(2, 278).Value = 288

# Copy the header style from the previous header cell (JP1) onto the new headers
$ws.Cells.Item(1, 276).Copy() | Out-Null
$ws.Range($ws.Cells.Item(1, 277), $ws.Cells.Item(1, 278)).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 277).Value = "Fatalities 12-07"
$ws.Cells.Item(1, 278).Value = "Fatalities 12-08"

# Match the column width formatting used across the rest of the sheet (12 characters)
$ws.Columns.Item(277).ColumnWidth = 11.17
$ws.Columns.Item(278).ColumnWidth = 11.17

# Update the active selection to match the saved workbook state
$ws.Range("C11").Select() | Out-Null
